$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166451692581177
$ws.Range("B1").Value = 2.429804086685181
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.371953010559082
$ws.Range("E1").Value = 1.235305190086365
